# correct typo in exam location (and incidental deck touch-ups)
#
# Summary of changes applied (see commit "correct typo in exam location"):
#   1. Slide 3:  "...5/5 in 54-340), ..."  ->  "...5/5 in 50-340), ..."
#   2. Slide 13: merge the "pencil-and-" / "paper" runs into one run of text
#      "pencil-and-paper" (no visible text change, just a run clean-up that
#      happens naturally when the range is reassigned).
#   3. Slide 14: touch up the "no puppies" run (no visible text change).
#   4. The slide master + all 11 slide layouts carry a cached
#      "datetimeFigureOut" field (Insert > Header & Footer date) that shows
#      the date the deck was last touched; bump it from 2/8/17 to 3/3/17.

$p = $ppt.ActivePresentation

function Update-DateField($shapes, $oldText, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# --- 1 & 2 & 3: text fixes on individual slides -----------------------

# Slide 3: "Grading: ... (3/3, 4/10, & 5/5 in 54-340), final exam 40%"
# fix the room-number typo 54-340 -> 50-340
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
$full3 = $tr3.Text
$idx3 = $full3.IndexOf("54-340")
if ($idx3 -ge 0) {
    $rng3 = $tr3.Characters($idx3 + 1, 6)
    $rng3.Text = "50-340"
}

# Slide 13: "... some computers vs. only pencil-and-paper"
# re-assign the run so the split "pencil-and-" / "paper" runs collapse
# into a single tidy run.
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$tr13 = $sh13.TextFrame.TextRange
$full13 = $tr13.Text
$idx13 = $full13.IndexOf("pencil-and-paper")
if ($idx13 -ge 0) {
    $rng13 = $tr13.Characters($idx13 + 1, 16)
    $rng13.Text = "pencil-and-paper"
}

# Slide 14: "... some puppy vs. no puppies"
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(2)
$tr14 = $sh14.TextFrame.TextRange
$full14 = $tr14.Text
$idx14 = $full14.IndexOf("no puppies")
if ($idx14 -ge 0) {
    $rng14 = $tr14.Characters($idx14 + 1, 10)
    $rng14.Text = "no puppies"
}

# --- 4: refresh the cached "datetimeFigureOut" footer date -------------
# on the slide master ...
$master = $p.SlideMaster
Update-DateField $master.Shapes "2/8/17" "3/3/17"

# ... and on every custom (slide) layout
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Update-DateField $layout.Shapes "2/8/17" "3/3/17"
}
